$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J8").Value = 0.2150495036779461
$ws.Range("I9").Value = 0.24
$ws.Range("H10").Value = 0.3087982760018804
$ws.Range("G11").Value = 0.32
$ws.Range("F12").Value = 0.4476495795507702
$ws.Range("E13").Value = 0.1088966743764388
$ws.Range("D14").Value = 0.1461563307127136
$ws.Range("C15").Value = 0.09547648014918764
$ws.Range("B16").Value = 0.0959495356205764
